# Actualización automática de tasas-transfi.xlsx

$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" message with new rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = @"
Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 12.66 = 51165.82 pesos
✅ 51165.82 pesos = 12.6 = 967.35 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%
"@

$newText = $newText.TrimEnd("`r", "`n")

$ws1.Range("A1").Value = $newText

# --- tasas: update the N10/O10/N12/O12 rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 79
$ws2.Range("O10").Value = 4042.1
$ws2.Range("N12").Value = 4060
$ws2.Range("O12").Value = 76.759
